$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 367.9091
$ws.Range("J53").Value = 850.3333
$ws.Range("L53").Value = 850.3333
$ws.Range("N53").Value = -2124.3333
$ws.Range("H86").Value = 7472.4443
$ws.Range("I86").Value = 11812.25
$ws.Range("J86").Value = 4000.6
$ws.Range("K86").Value = 11812.25
$ws.Range("L86").Value = 4000.6
$ws.Range("M86").Value = -10689.25
$ws.Range("N86").Value = -6246.6
$ws.Range("H89").Value = 7472.4443
$ws.Range("I89").Value = 11812.25
$ws.Range("J89").Value = 4000.6
$ws.Range("K89").Value = 59061.25
$ws.Range("L89").Value = 20003
$ws.Range("M89").Value = -53445.25
$ws.Range("N89").Value = -31235
$ws.Range("H141").Value = 1432.4
$ws.Range("I141").Value = 1432.4
$ws.Range("K141").Value = 4297.200000000001
$ws.Range("M141").Value = 882.7999999999993

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7734.6665
$ws.Range("I45").Value = 9028.286
$ws.Range("K45").Value = 9028.286
$ws.Range("M45").Value = -8651.286
$ws.Range("H61").Value = 55557316
$ws.Range("I61").Value = 66668380
$ws.Range("K61").Value = 66668380
$ws.Range("M61").Value = -66668168
$ws.Range("H74").Value = 27780692
$ws.Range("I74").Value = 28574412
$ws.Range("K74").Value = 28574412
$ws.Range("M74").Value = -28573538
$ws.Range("H77").Value = 27780692
$ws.Range("I77").Value = 28574412
$ws.Range("K77").Value = 142872060
$ws.Range("M77").Value = -142867692
$ws.Range("H110").Value = 201246.6
$ws.Range("I110").Value = 201246.6
$ws.Range("K110").Value = 201246.6
$ws.Range("M110").Value = -199201.6
$ws.Range("H136").Value = 55557316
$ws.Range("I136").Value = 66668380
$ws.Range("K136").Value = 200005140
$ws.Range("M136").Value = -200002590

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1866
$ws.Range("I22").Value = 1244
$ws.Range("J22").Value = 2073.3333
$ws.Range("K22").Value = 1244
$ws.Range("L22").Value = 2073.3333
$ws.Range("M22").Value = -1071
$ws.Range("N22").Value = -2419.3333
$ws.Range("H99").Value = 1533.1111
$ws.Range("I99").Value = 1493.6957
$ws.Range("J99").Value = 1759.75
$ws.Range("K99").Value = 1493.6957
$ws.Range("L99").Value = 1759.75
$ws.Range("M99").Value = 4.304300000000012
$ws.Range("N99").Value = -4755.75
$ws.Range("H134").Value = 33334492
$ws.Range("I134").Value = 33334492
$ws.Range("K134").Value = 100003476
$ws.Range("M134").Value = -100000941

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 176999.5
$ws.Range("I6").Value = 234333
$ws.Range("K6").Value = 234333
$ws.Range("M6").Value = -234220
$ws.Range("H31").Value = 4318
$ws.Range("I31").Value = 5782.4
$ws.Range("K31").Value = 5782.4
$ws.Range("M31").Value = -5487.4
$ws.Range("H34").Value = 4318
$ws.Range("I34").Value = 5782.4
$ws.Range("K34").Value = 5782.4
$ws.Range("M34").Value = -5580.4
$ws.Range("H58").Value = 62514704
$ws.Range("I58").Value = 83351770
$ws.Range("K58").Value = 83351770
$ws.Range("M58").Value = -83351567
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H86").Value = 14933.333
$ws.Range("I86").Value = 12500
$ws.Range("J86").Value = 15237.5
$ws.Range("K86").Value = 12500
$ws.Range("L86").Value = 15237.5
$ws.Range("M86").Value = -11377
$ws.Range("N86").Value = -17483.5
$ws.Range("H89").Value = 14933.333
$ws.Range("I89").Value = 12500
$ws.Range("J89").Value = 15237.5
$ws.Range("K89").Value = 62500
$ws.Range("L89").Value = 76187.5
$ws.Range("M89").Value = -56884
$ws.Range("N89").Value = -87419.5
$ws.Range("H99").Value = 14658.667
$ws.Range("J99").Value = 6989
$ws.Range("L99").Value = 6989
$ws.Range("N99").Value = -9985
$ws.Range("H126").Value = 14658.667
$ws.Range("J126").Value = 6989
$ws.Range("L126").Value = 20967
$ws.Range("N126").Value = -25907
$ws.Range("H136").Value = 62514704
$ws.Range("I136").Value = 83351770
$ws.Range("K136").Value = 250055310
$ws.Range("M136").Value = -250052760

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 120.42857
$ws.Range("I23").Value = 109.333336
$ws.Range("J23").Value = 128.75
$ws.Range("K23").Value = 328.000008
$ws.Range("L23").Value = 386.25
$ws.Range("M23").Value = -93.00000799999998
$ws.Range("N23").Value = -856.25
$ws.Range("H75").Value = 300.2857
$ws.Range("I75").Value = 298.6
$ws.Range("J75").Value = 304.5
$ws.Range("K75").Value = 895.8000000000001
$ws.Range("L75").Value = 913.5
$ws.Range("M75").Value = 102.1999999999999
$ws.Range("N75").Value = -2909.5
$ws.Range("H78").Value = 300.2857
$ws.Range("I78").Value = 298.6
$ws.Range("J78").Value = 304.5
$ws.Range("K78").Value = 2687.4
$ws.Range("L78").Value = 2740.5
$ws.Range("M78").Value = 2304.6
$ws.Range("N78").Value = -12724.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 106.9375
$ws.Range("I2").Value = 46.88889
$ws.Range("K2").Value = 46.88889
$ws.Range("M2").Value = 66.11111
$ws.Range("H20").Value = 30000
$ws.Range("J20").Value = 30000
$ws.Range("L20").Value = 30000
$ws.Range("N20").Value = -30490
$ws.Range("H70").Value = 4389
$ws.Range("I70").Value = 4389
$ws.Range("K70").Value = 4389
$ws.Range("M70").Value = -4119
$ws.Range("H73").Value = 4389
$ws.Range("I73").Value = 4389
$ws.Range("K73").Value = 4389
$ws.Range("M73").Value = -3453
$ws.Range("H97").Value = 1732.9048
$ws.Range("I97").Value = 1636.4166
$ws.Range("K97").Value = 1636.4166
$ws.Range("M97").Value = -1140.4166

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 27999
$ws.Range("I20").Value = 27999
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 27999
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -27773
$ws.Range("N20").ClearContents()
$ws.Range("H21").Value = 1700
$ws.Range("I21").Value = 1700
$ws.Range("K21").Value = 1700
$ws.Range("M21").Value = -1526
$ws.Range("H22").Value = 2897.5454
$ws.Range("I22").Value = 3110.4443
$ws.Range("J22").Value = 1939.5
$ws.Range("K22").Value = 3110.4443
$ws.Range("L22").Value = 1939.5
$ws.Range("M22").Value = -2815.4443
$ws.Range("N22").Value = -2529.5
$ws.Range("H27").Value = 2897.5454
$ws.Range("I27").Value = 3110.4443
$ws.Range("J27").Value = 1939.5
$ws.Range("K27").Value = 3110.4443
$ws.Range("L27").Value = 1939.5
$ws.Range("M27").Value = -3003.4443
$ws.Range("N27").Value = -2153.5
$ws.Range("H46").Value = 2255
$ws.Range("I46").Value = 2272.2222
$ws.Range("J46").Value = 2100
$ws.Range("K46").Value = 2272.2222
$ws.Range("L46").Value = 2100
$ws.Range("M46").Value = -2084.2222
$ws.Range("N46").Value = -2476
$ws.Range("H94").Value = 135000
$ws.Range("J94").Value = 135000
$ws.Range("L94").Value = 135000
$ws.Range("N94").Value = -136352

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 6229.8
$ws.Range("I32").Value = 7537.25
$ws.Range("J32").Value = 1000
$ws.Range("K32").Value = 7537.25
$ws.Range("L32").Value = 1000
$ws.Range("M32").Value = -7220.25
$ws.Range("N32").Value = -1634
$ws.Range("H100").Value = 2119.75
$ws.Range("I100").Value = 2119.75
$ws.Range("K100").Value = 4239.5
$ws.Range("M100").Value = -3698.5
$ws.Range("H122").Value = 2222.2083
$ws.Range("I122").Value = 1833.8125
$ws.Range("K122").Value = 5501.4375
$ws.Range("M122").Value = -3051.4375
